$d = $word.ActiveDocument

function Split-RunsInRange($doc, $baseStart, $spans) {
    foreach ($s in $spans) {
        $a = $s[0]
        $b = $s[1]
        $pt = $doc.Range($baseStart + $a, $baseStart + $b)
        $pt.Bold = 1
        $pt.Bold = 0
    }
}

# ---------------------------------------------------------------------------
# 1) Author line: "L.C. Walker, L.A. Rogers, S.C. Anderson, and D.R. Haggarty"
#    becomes "Leah C. Walker, Luke A. Rogers, Sean C. Anderson, and Dana R. Haggarty"
#    split across 10 runs.
# ---------------------------------------------------------------------------
$authorOld = "L.C. Walker, L.A. Rogers, S.C. Anderson, and D.R. Haggarty"
$authorNew = "Leah C. Walker, Luke A. Rogers, Sean C. Anderson, and Dana R. Haggarty"

$rngAuthor = $d.Content
$rngAuthor.Find.Execute($authorOld, $true, $false, $false, $false, $false, $true, 1, $false, $authorNew, 2) | Out-Null

$rngAuthorFound = $d.Content
$rngAuthorFound.Find.Execute($authorNew, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$authorBase = $rngAuthorFound.Start

$authorSpans = @(
    @(0, 1),    # L
    @(1, 5),    # eah 
    @(5, 14),   # C. Walker
    @(14, 17),  # , L
    @(17, 21),  # uke 
    @(21, 33),  # A. Rogers, S
    @(33, 37),  # ean 
    @(37, 55),  # C. Anderson, and D
    @(55, 59),  # ana 
    @(59, 70)   # R. Haggarty
)
Split-RunsInRange $d $authorBase $authorSpans

# ---------------------------------------------------------------------------
# 2) Address block:
#    "Nanaimo, British Columbia" (own paragraph) -> "3190 Hammond Bay Road"
#    "V9T 6N7" + trailing " " run (own paragraph) -> "Nanaimo, British Columbia" / ", " / "V9T 6N7"
# ---------------------------------------------------------------------------
$rngCity = $d.Content
$rngCity.Find.Execute("Nanaimo, British Columbia", $true, $false, $false, $false, $false, $true, 1, $false, "3190 Hammond Bay Road", 2) | Out-Null

# Merge "V9T 6N7" and the trailing space run into a single run without the space.
$rngPostal = $d.Content
$rngPostal.Find.Execute("V9T 6N7 ", $true, $false, $false, $false, $false, $true, 1, $false, "V9T 6N7", 2) | Out-Null

# Prepend the city/province text in front of the postal code.
$addressNew = "Nanaimo, British Columbia, V9T 6N7"
$rngAddress = $d.Content
$rngAddress.Find.Execute("V9T 6N7", $true, $false, $false, $false, $false, $true, 1, $false, $addressNew, 2) | Out-Null

$rngAddressFound = $d.Content
$rngAddressFound.Find.Execute($addressNew, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$addressBase = $rngAddressFound.Start

$addressSpans = @(
    @(0, 25),   # Nanaimo, British Columbia
    @(25, 27),  # ", "
    @(27, 34)   # V9T 6N7
)
Split-RunsInRange $d $addressBase $addressSpans

Write-Output "Edit complete"
